$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")


# --- Reset target range formatting baseline (General numfmt, not bold, no borders, default align) ---
$fullRange = $ws.Range("E4:L28")
$fullRange.Font.Bold = $false
$fullRange.HorizontalAlignment = -4131  # xlGeneral
$fullRange.NumberFormat = "General"
$none = -4142  # xlLineStyleNone
$fullRange.Borders.Item(7).LineStyle = $none
$fullRange.Borders.Item(8).LineStyle = $none
$fullRange.Borders.Item(9).LineStyle = $none
$fullRange.Borders.Item(10).LineStyle = $none

# --- Set cell values ---
$ws.Range("E4").Value = "Register Convention--Logo"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = "North"
$ws.Range("H6").Value = "Red"
$ws.Range("I6").Value = 16
$ws.Range("J6").Value = "curr. Direction [0-3]"
$ws.Range("L6").Value = "Previous"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "East"
$ws.Range("H7").Value = "Green"
$ws.Range("I7").Value = 17
$ws.Range("J7").Value = "line color"
$ws.Range("L7").Value = "Previous"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = "South"
$ws.Range("H8").Value = "Blue"
$ws.Range("I8").Value = 18
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = "West"
$ws.Range("H9").Value = "White"
$ws.Range("I9").Value = 19
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = "arg0"
$ws.Range("I10").Value = 20
$ws.Range("J10").Value = "tmp4"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = "arg1"
$ws.Range("I11").Value = 21
$ws.Range("J11").Value = "tmp5"
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = "tmp0"
$ws.Range("I12").Value = 22
$ws.Range("J12").Value = "tmp6"
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "tmp1"
$ws.Range("I13").Value = 23
$ws.Range("J13").Value = "tmp7"
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = "tmp2"
$ws.Range("I14").Value = 24
$ws.Range("J14").Value = "tmp8"
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = "tmp3"
$ws.Range("I15").Value = 25
$ws.Range("J15").Value = "tmp9"
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = "x-coord"
$ws.Range("H16").Value = "Current"
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = "tmp10"
$ws.Range("E17").Value = 11
$ws.Range("F17").Value = "y-coord"
$ws.Range("H17").Value = "Current"
$ws.Range("I17").Value = 27
$ws.Range("J17").Value = "tmp11"
$ws.Range("E18").Value = 12
$ws.Range("F18").Value = "curr. Direction [0-3]"
$ws.Range("H18").Value = "Current"
$ws.Range("I18").Value = 28
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = "line color"
$ws.Range("H19").Value = "Current"
$ws.Range("I19").Value = 29
$ws.Range("J19").Value = "DMEM stack ptr"
$ws.Range("E20").Value = 14
$ws.Range("F20").Value = "x-coord"
$ws.Range("H20").Value = "Previous"
$ws.Range("I20").Value = 30
$ws.Range("J20").Value = "pen down/up"
$ws.Range("E21").Value = 15
$ws.Range("F21").Value = "y-coord"
$ws.Range("H21").Value = "Previous"
$ws.Range("I21").Value = 31
$ws.Range("J21").Value = "JR"

# --- Group cells by formatting signature and apply ---

# cells: E4
foreach ($ref in @("E4")) {
  $rg = $ws.Range($ref)
  $rg.Font.Bold = $true
}

# cells: E6, I6
foreach ($ref in @("E6", "I6")) {
  $rg = $ws.Range($ref)
  $rg.NumberFormat = "`"$`"#,##0_);[Red](`"$`"#,##0)"
  $rg.Borders.Item(7).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(7).ColorIndex = 1
  $rg.Borders.Item(8).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(8).ColorIndex = 1
}

# cells: F6
foreach ($ref in @("F6")) {
  $rg = $ws.Range($ref)
  $rg.Font.Bold = $true
  $rg.HorizontalAlignment = -4108  # xlCenter
  $rg.Borders.Item(8).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(8).ColorIndex = 1
}

# cells: G6
foreach ($ref in @("G6")) {
  $rg = $ws.Range($ref)
  $rg.HorizontalAlignment = -4108  # xlCenter
  $rg.Borders.Item(8).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(8).ColorIndex = 1
}

# cells: H6, L6
foreach ($ref in @("H6", "L6")) {
  $rg = $ws.Range($ref)
  $rg.HorizontalAlignment = -4108  # xlCenter
  $rg.Borders.Item(8).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(8).ColorIndex = 1
  $rg.Borders.Item(10).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(10).ColorIndex = 1
}

# cells: J6, K6
foreach ($ref in @("J6", "K6")) {
  $rg = $ws.Range($ref)
  $rg.Borders.Item(8).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(8).ColorIndex = 1
}

# cells: E7, I7, E8, I8, E9, I9, E10, I10, E11, I11, E12, I12, E13, I13, E14, I14, E15, I15, E16, I16, E17, I17, E18, I18, E19, I19, E20, I20
foreach ($ref in @("E7", "I7", "E8", "I8", "E9", "I9", "E10", "I10", "E11", "I11", "E12", "I12", "E13", "I13", "E14", "I14", "E15", "I15", "E16", "I16", "E17", "I17", "E18", "I18", "E19", "I19", "E20", "I20")) {
  $rg = $ws.Range($ref)
  $rg.NumberFormat = "`"$`"#,##0_);[Red](`"$`"#,##0)"
  $rg.Borders.Item(7).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(7).ColorIndex = 1
}

# cells: F7, F8, F9
foreach ($ref in @("F7", "F8", "F9")) {
  $rg = $ws.Range($ref)
  $rg.Font.Bold = $true
  $rg.HorizontalAlignment = -4108  # xlCenter
}

# cells: G7, G8, G9, G10
foreach ($ref in @("G7", "G8", "G9", "G10")) {
  $rg = $ws.Range($ref)
  $rg.HorizontalAlignment = -4108  # xlCenter
}

# cells: H7, L7, H8, H9, H10, H16, H17, H18, H19, H20
foreach ($ref in @("H7", "L7", "H8", "H9", "H10", "H16", "H17", "H18", "H19", "H20")) {
  $rg = $ws.Range($ref)
  $rg.HorizontalAlignment = -4108  # xlCenter
  $rg.Borders.Item(10).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(10).ColorIndex = 1
}

# cells: J7, K7, J8, K8, J9, K9, F10, J10, K10, F11, G11, J11, K11, F12, G12, J12, K12, F13, G13, J13, K13, F14, G14, J14, K14, F15, G15, J15, K15, F16, G16, J16, K16, F17, G17, J17, K17, F18, G18, J18, K18, F19, G19, J19, K19, F20, G20, J20, K20
foreach ($ref in @("J7", "K7", "J8", "K8", "J9", "K9", "F10", "J10", "K10", "F11", "G11", "J11", "K11", "F12", "G12", "J12", "K12", "F13", "G13", "J13", "K13", "F14", "G14", "J14", "K14", "F15", "G15", "J15", "K15", "F16", "G16", "J16", "K16", "F17", "G17", "J17", "K17", "F18", "G18", "J18", "K18", "F19", "G19", "J19", "K19", "F20", "G20", "J20", "K20")) {
  $rg = $ws.Range($ref)
}

# cells: L8, L9, L10, H11, L11, H12, L12, H13, L13, H14, L14, H15, L15, L16, L17, L18, L19
foreach ($ref in @("L8", "L9", "L10", "H11", "L11", "H12", "L12", "H13", "L13", "H14", "L14", "H15", "L15", "L16", "L17", "L18", "L19")) {
  $rg = $ws.Range($ref)
  $rg.Borders.Item(10).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(10).ColorIndex = 1
}

# cells: L20
foreach ($ref in @("L20")) {
  $rg = $ws.Range($ref)
  $rg.NumberFormat = "mmm-yy"
  $rg.Borders.Item(10).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(10).ColorIndex = 1
}

# cells: E21, I21
foreach ($ref in @("E21", "I21")) {
  $rg = $ws.Range($ref)
  $rg.NumberFormat = "`"$`"#,##0_);[Red](`"$`"#,##0)"
  $rg.Borders.Item(7).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(7).ColorIndex = 1
  $rg.Borders.Item(9).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(9).ColorIndex = 1
}

# cells: F21, G21, J21, K21
foreach ($ref in @("F21", "G21", "J21", "K21")) {
  $rg = $ws.Range($ref)
  $rg.Borders.Item(9).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(9).ColorIndex = 1
}

# cells: H21
foreach ($ref in @("H21")) {
  $rg = $ws.Range($ref)
  $rg.HorizontalAlignment = -4108  # xlCenter
  $rg.Borders.Item(9).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(9).ColorIndex = 1
  $rg.Borders.Item(10).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(10).ColorIndex = 1
}

# cells: L21
foreach ($ref in @("L21")) {
  $rg = $ws.Range($ref)
  $rg.Borders.Item(9).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(9).ColorIndex = 1
  $rg.Borders.Item(10).LineStyle = -4119  # xlDouble
  $rg.Borders.Item(10).ColorIndex = 1
}

# cells: E23, E24, E25, E26, E27, E28
foreach ($ref in @("E23", "E24", "E25", "E26", "E27", "E28")) {
  $rg = $ws.Range($ref)
  $rg.NumberFormat = "`"$`"#,##0_);[Red](`"$`"#,##0)"
}

# --- View / selection state ---
$ws.Range("N11").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 4
